$d = $word.ActiveDocument

# Paragraph 1 currently holds the old "waiting for call" rule text. It is
# replaced with new guidance that only counts assistant messages, plus a
# quoted status phrase, and a brand-new paragraph with the follow-up
# confirmation rule is inserted right after it.

$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Анализируй только сообщения, отправленные assistant, для определения подтверждения звонка. Если подтверждение отсутствует, а клиент только предлагает время или спрашивает уточнения, выводи статус «разговор продолжается»"

# Insert a brand-new paragraph right after it for the assistant-confirmation rule.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Если assistant отправляет сообщение, в котором он подтверждает выбранное время звонка, то напиши «статус ожидает звонка»"
